$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column J (this pushes the old I column -- the
#    total_fertN_recovery / (H+F) percent formula -- into J, and old J
#    (Year) into K). We'll then fill the now-blank I with a brand new
#    "total_fertN_recovery" = G+E formula, and relabel the headers.
# ---------------------------------------------------------------------------
$ws.Columns("J:J").Insert()

# ---------------------------------------------------------------------------
# 2. New column I formulas: total_fertN_recovery = crop_recovery_fert_N (G)
#    + soil_recovery_fert_N (E)
# ---------------------------------------------------------------------------
$ws.Range("I2").Formula = "=G2+E2"
$ws.Range("I3:I25").Formula = "=G3+E3"

# ---------------------------------------------------------------------------
# 3. Column J keeps the original total formula (H+F), now labelled as the
#    percent version.
# ---------------------------------------------------------------------------
$ws.Range("J2").Formula = "=H2+F2"
$ws.Range("J3:J25").Formula = "=H3+F3"

# ---------------------------------------------------------------------------
# 4. Header row relabeling (shared-string append order matters: this must
#    happen before the Stage/Maturity column so "total_fertN_recovery_percent"
#    lands before "Stage"/"Maturity" in the shared-string table, matching the
#    canonical save order).
# ---------------------------------------------------------------------------
$ws.Range("I1").Value = "total_fertN_recovery"
$ws.Range("J1").Value = "total_fertN_recovery_percent"
$ws.Range("K1").Value = "Year"

# ---------------------------------------------------------------------------
# 5. Column D: "DAS" / "127" -> "Stage" / "Maturity"
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Stage"
$ws.Range("D2:D25").Value = "Maturity"

# ---------------------------------------------------------------------------
# 6. New blank column M, formatted with the same "0.0" number style used by
#    the rest of the numeric columns (style index 1 -- numFmt 164).
# ---------------------------------------------------------------------------
$ws.Range("M2:M25").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 7. Column widths: widen I, add a width for the new J column.
# ---------------------------------------------------------------------------
$ws.Columns("I:I").ColumnWidth = 21.7
$ws.Columns("J:J").ColumnWidth = 28.7

# ---------------------------------------------------------------------------
# 8. Selection moves to D2:D25 (active cell D2)
# ---------------------------------------------------------------------------
$ws.Range("D2:D25").Select()
